$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 onto the two new header cells so they
# reuse the existing bold/centered/bordered style instead of minting a new one.
$ws.Range("H1:H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I (I0) and J (IF) columns, rows 2-12
$data = @(
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(7, 8),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(2, 3),
    @(8, 8),
    @(8, 8),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
